# Add a new "2023" column (column T) to the Hepatitis B incidence sheet,
# mirroring the existing "2022" column (S) for layout/formatting, and
# bump row 4's height to match the taller header row used once the new
# column data pushed the sheet slightly differently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 gets an explicit custom height in the edited workbook.
$ws.Rows.Item(4).RowHeight = 16.5

# Clone the formatting of column S (years 2007-2022) onto the new column T
# (year 2023) before writing values, so number formats / fonts / borders
# match the rest of the table.
$ws.Range("S3:S33").Copy()
$ws.Range("T3:T33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header year
$ws.Cells.Item(3, 20).Value = 2023

# Data values for the new 2023 column, row by row.
$ws.Cells.Item(4, 20).Value = 2.3381104968484805
$ws.Cells.Item(5, 20).Value = 2.0344672190198714
$ws.Cells.Item(6, 20).Value = 2.6483752218014245
$ws.Cells.Item(7, 20).Value = 3.9852372948902328
$ws.Cells.Item(8, 20).Value = 4.5532396299967433
$ws.Cells.Item(9, 20).Value = 3.4291318466903733
$ws.Cells.Item(10, 20).Value = 1.2089851778417198
$ws.Cells.Item(11, 20).Value = 1.521116134174612
$ws.Cells.Item(12, 20).Value = 0.9008846687447073
$ws.Cells.Item(13, 20).Value = 3.694303753043183
$ws.Cells.Item(14, 20).Value = 4.0607488020791038
$ws.Cells.Item(15, 20).Value = 3.327319511401615
$ws.Cells.Item(16, 20).Value = 0.32236434908190637
$ws.Cells.Item(17, 20).Value = 0
$ws.Cells.Item(18, 20).Value = 0.63756806039044667
$ws.Cells.Item(19, 20).Value = 2.1691385808410835
$ws.Cells.Item(20, 20).Value = 1.5024572004578396
$ws.Cells.Item(21, 20).Value = 2.8259763748375066
$ws.Cells.Item(22, 20).Value = 6.1744985943935555
$ws.Cells.Item(23, 20).Value = 4.3993752887090034
$ws.Cells.Item(24, 20).Value = 7.9169155696940479
$ws.Cells.Item(25, 20).Value = 2.8763040791558883
$ws.Cells.Item(26, 20).Value = 1.4751329463567904
$ws.Cells.Item(27, 20).Value = 4.2954684675262591
$ws.Cells.Item(28, 20).Value = 1.8177568880002077
$ws.Cells.Item(29, 20).Value = 1.581380197008345
$ws.Cells.Item(30, 20).Value = 2.103608453446189
$ws.Cells.Item(31, 20).Value = 1.3736037318066185
$ws.Cells.Item(32, 20).Value = 2.249820014398848
$ws.Cells.Item(33, 20).Value = 0.53701655085009725
